$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (subject/column identifiers) for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - updated meanEMG legmaxROM values for columns B:E
$ws.Range("B2").Value = 17.312042658875065
$ws.Range("C2").Value = 44.663834416666418
$ws.Range("D2").Value = 51.934247515874574
$ws.Range("E2").Value = 48.87056625032141

# Row 3 (STR) - updated meanEMG legmaxROM values for columns B:E
$ws.Range("B3").Value = 30.333193684649491
$ws.Range("C3").Value = 68.291712523665268
$ws.Range("D3").Value = 69.127959583571368
$ws.Range("E3").Value = 47.032783547921092

# Update the active selection on the sheet to match the new data range of interest
$ws.Range("B1:E3").Select() | Out-Null
